# Secdep Loan, Saving, RD scenarios
# Updates the "Summary" and "Repayment schedule" sheets of the bulk JLG
# loan workbook, and moves the active tab / selections accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: principal due / outstanding for the first instalment
# row dropped from 275.05 to 271.52 (Original + Outstanding columns).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A3").Value = 271.52
$summary.Range("E3").Value = 271.52
# Move this sheet's own cursor (it is not the active tab any more).
$summary.Range("D6").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet: instalments 10-12 recomputed.
# ---------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Repayment schedule")

# Instalment 10 (row 12)
$sched.Range("B12").Value = 14
$sched.Range("C12").Value = 42157
$sched.Range("F12").Value = 762.06
$sched.Range("G12").Value = 1534.62
$sched.Range("H12").Value = 10.57

# Instalment 11 (row 13)
$sched.Range("B13").Value = 14
$sched.Range("F13").Value = 765.57
$sched.Range("G13").Value = 769.05
$sched.Range("H13").Value = 7.06

# Instalment 12 (row 14)
$sched.Range("F14").Value = 769.05
$sched.Range("H14").Value = 3.54
$sched.Range("K14").Value = 772.59
$sched.Range("Q14").Value = 772.59

# Make "Repayment schedule" the active sheet/tab, replacing
# "BulkJLGLoanInput", and set its new selection.
$sched.Activate()
$sched.Range("J16").Select()
